$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'302.97"
$ws.Range("E2").Value = "'0.09%"
$ws.Range("D3").Value = "'32.67"
$ws.Range("E3").Value = "'1.69%"
$ws.Range("D4").Value = "'5.079"
$ws.Range("E4").Value = "'-1.06%"
$ws.Range("D5").Value = "'0.07717"
$ws.Range("E5").Value = "'-1.51%"
$ws.Range("D6").Value = "'2.085"
$ws.Range("E6").Value = "'-8.33%"
$ws.Range("D7").Value = "'7.905"
$ws.Range("E7").Value = "'1.07%"
$ws.Range("B8").Value = "MXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D8").Value = "'0.9229"
$ws.Range("E8").Value = "'-0.49%"
$ws.Range("B9").Value = "WazirX"
$ws.Range("C9").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D9").Value = "'0.1767"
$ws.Range("E9").Value = "'0.10%"
$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D10").Value = "'0.07977"
$ws.Range("E10").Value = "'3.43%"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "'0.08584"
$ws.Range("E11").Value = "'-3.32%"
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").Value = "'0.03062"
$ws.Range("E12").Value = "'0.41%"
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").Value = "'0.09988"
$ws.Range("E13").Value = "'-0.09%"
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D14").Value = "'0.001516"
$ws.Range("E14").Value = "'0.39%"
$ws.Range("B15").Value = "TigerCash"
$ws.Range("C15").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D15").Value = "'0.005675"
$ws.Range("E15").Value = "'-5.77%"
$ws.Range("B16").Value = "UpBots"
$ws.Range("C16").Value = "https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt"
$ws.Range("D16").Value = "'0.007498"
$ws.Range("E16").Value = "'2,116.77%"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").Value = "'3.469"
$ws.Range("E17").Value = "'0.14%"
$ws.Range("B18").Value = "GateToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D18").Value = "'3.798"
$ws.Range("E18").Value = "'-0.15%"
$ws.Range("E19").Value = "'-4.28%"
$ws.Range("D20").Value = "'0.3352"
$ws.Range("E20").Value = "'2.44%"
$ws.Range("D21").Value = "'0.1301"
$ws.Range("E21").Value = "'-3.53%"
$ws.Range("D22").Value = "'4.394"
$ws.Range("E22").Value = "'2.93%"
$ws.Range("D23").Value = "'0.1978"
$ws.Range("E23").Value = "'9.90%"
$ws.Range("D24").Value = "'0.04536"
$ws.Range("E24").Value = "'-1.17%"
$ws.Range("D25").Value = "'0.001234"
$ws.Range("E25").Value = "'-1.75%"
$ws.Range("D26").Value = "'0.004145"
$ws.Range("E26").Value = "'-7.91%"
$ws.Range("D27").Value = "'0.0001251"
$ws.Range("E27").Value = "'-0.08%"
$ws.Range("D39").Value = "'0.01726"
$ws.Range("E39").Value = "'-3.01%"
$ws.Range("D40").Value = "'0.04719"
$ws.Range("E40").Value = "'-0.09%"
$ws.Range("D41").Value = "'0.007482"
$ws.Range("E41").Value = "'3.30%"
$ws.Range("D42").Value = "'0.1361"
$ws.Range("E42").Value = "'-0.58%"
$ws.Range("D43").Value = "'0.002332"
$ws.Range("E43").Value = "'9.70%"
$ws.Range("D44").Value = "'0.01061"
$ws.Range("E44").Value = "'-4.25%"
$ws.Range("D45").Value = "'0.00006207"
$ws.Range("E45").Value = "'-0.34%"
$ws.Range("E46").Value = "'-0.10%"
$ws.Range("D47").Value = "'1.101"
$ws.Range("E47").Value = "'47.15%"
$ws.Range("D48").Value = "'0.003004"
$ws.Range("E48").Value = "'-6.24%"
$ws.Range("D49").Value = "'0.00002101"
$ws.Range("E49").Value = "'-0.10%"
$ws.Range("D50").Value = "'0.0002001"
$ws.Range("E50").Value = "'-0.10%"
